$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -831
$ws.Range("N13").ClearContents()

# Row 15
$ws.Range("H15").Value = 122.76
$ws.Range("I15").Value = 122.76
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 368.28
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -199.28

# Row 117
$ws.Range("H117").Value = 43242
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 43242
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 43242
$ws.Range("N117").Value = -52420

# Row 132
$ws.Range("H132").Value = 6582605
$ws.Range("I132").Value = 3980.5356
$ws.Range("J132").Value = 25002754
$ws.Range("K132").Value = 11941.6068
$ws.Range("L132").Value = 75008262
$ws.Range("M132").Value = -9411.606800000001
$ws.Range("N132").Value = -75013322

# Row 137
$ws.Range("H137").Value = 1068.7273
$ws.Range("I137").Value = 1229.9
$ws.Range("J137").Value = 934.4167
$ws.Range("K137").Value = 3689.7
$ws.Range("L137").Value = 2803.2501
$ws.Range("M137").Value = -1139.7

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 254.25
$ws.Range("I5").Value = 231.90909
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 231.90909
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -119.90909

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()

# Row 32
$ws.Range("H32").Value = 3877.44
$ws.Range("I32").Value = 3584.6316
$ws.Range("J32").Value = 9440.799999999999
$ws.Range("K32").Value = 3584.6316
$ws.Range("L32").Value = 9440.799999999999
$ws.Range("M32").Value = -3297.6316
$ws.Range("N32").Value = -10014.8

# Row 61
$ws.Range("H61").Value = 2910.92
$ws.Range("I61").Value = 2989.9583
$ws.Range("J61").Value = 1014
$ws.Range("K61").Value = 2989.9583
$ws.Range("L61").Value = 1014
$ws.Range("M61").Value = -2777.9583

# Row 102
$ws.Range("H102").Value = 1098.1666
$ws.Range("I102").Value = 1098.1666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1098.1666
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 523.8334

# Row 132
$ws.Range("H132").Value = 7145083.5
$ws.Range("I132").Value = 13891053
$ws.Range("J132").Value = 2292.7058
$ws.Range("K132").Value = 41673159
$ws.Range("L132").Value = 6878.117400000001
$ws.Range("M132").Value = -41670629
$ws.Range("N132").Value = -11938.1174

# Row 136
$ws.Range("H136").Value = 2910.92
$ws.Range("I136").Value = 2989.9583
$ws.Range("J136").Value = 1014
$ws.Range("K136").Value = 8969.874899999999
$ws.Range("L136").Value = 3042
$ws.Range("M136").Value = -6419.874899999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 254.25
$ws.Range("I4").Value = 231.90909
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 231.90909
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -116.90909

# Row 23
$ws.Range("H23").Value = 11000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 11000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 11000
$ws.Range("N23").Value = -11566

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 365.375
$ws.Range("I5").Value = 261.4
$ws.Range("J5").Value = 538.6667
$ws.Range("K5").Value = 261.4
$ws.Range("L5").Value = 538.6667
$ws.Range("M5").Value = -149.4
$ws.Range("N5").Value = -762.6667

# Row 8
$ws.Range("H8").Value = 634.25
$ws.Range("I8").Value = 509
$ws.Range("J8").Value = 1010
$ws.Range("K8").Value = 509
$ws.Range("L8").Value = 1010
$ws.Range("M8").Value = -369
$ws.Range("N8").Value = -1290

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()

# Row 14
$ws.Range("H14").Value = 4666.6665
$ws.Range("I14").Value = 6000
$ws.Range("J14").Value = 4000
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = -5830
$ws.Range("N14").Value = -4340

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 19
$ws.Range("H19").Value = 107
$ws.Range("I19").Value = 107
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 107
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 63
$ws.Range("N19").ClearContents()

# Row 24
$ws.Range("H24").Value = 107
$ws.Range("I24").Value = 107
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 107
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 63
$ws.Range("N24").ClearContents()

# Row 31
$ws.Range("H31").Value = 2402.3948
$ws.Range("I31").Value = 1892.159
$ws.Range("J31").Value = 3103.9688
$ws.Range("K31").Value = 1892.159
$ws.Range("L31").Value = 3103.9688
$ws.Range("M31").Value = -1597.159
$ws.Range("N31").Value = -3693.9688

# Row 34
$ws.Range("H34").Value = 2402.3948
$ws.Range("I34").Value = 1892.159
$ws.Range("J34").Value = 3103.9688
$ws.Range("K34").Value = 1892.159
$ws.Range("L34").Value = 3103.9688
$ws.Range("M34").Value = -1690.159
$ws.Range("N34").Value = -3507.9688

# Row 94
$ws.Range("H94").Value = 6149
$ws.Range("I94").Value = 4202
$ws.Range("J94").Value = 7365.875
$ws.Range("K94").Value = 4202
$ws.Range("L94").Value = 7365.875
$ws.Range("M94").Value = -3751
$ws.Range("N94").Value = -8267.875

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 2451.5
$ws.Range("I19").Value = 2901
$ws.Range("J19").Value = 2002
$ws.Range("K19").Value = 8703
$ws.Range("L19").Value = 6006
$ws.Range("M19").Value = -8529
$ws.Range("N19").Value = -6354

# Row 56
$ws.Range("H56").Value = 4971.4243
$ws.Range("I56").Value = 4971.4243
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4971.4243
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -4441.4243

# Row 122
$ws.Range("H122").Value = 1301.5769
$ws.Range("I122").Value = 1116.0588
$ws.Range("J122").Value = 1652
$ws.Range("K122").Value = 10044.5292
$ws.Range("L122").Value = 14868
$ws.Range("M122").Value = -7594.529200000001
$ws.Range("N122").Value = -19768

# Row 125
$ws.Range("H125").Value = 5274.9375
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 5274.9375
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 15824.8125
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -25664.8125

# Row 131
$ws.Range("H131").Value = 2004.2609
$ws.Range("I131").Value = 2851.0588
$ws.Range("J131").Value = 1507.862
$ws.Range("K131").Value = 8553.1764
$ws.Range("L131").Value = 4523.586
$ws.Range("M131").Value = -3513.1764
$ws.Range("N131").Value = -14603.586

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 8501.75
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 8501.75
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 8501.75
$ws.Range("N21").Value = -8847.75

# Row 30
$ws.Range("H30").Value = 8501.75
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 8501.75
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 8501.75
$ws.Range("N30").Value = -8711.75

# Row 97
$ws.Range("H97").Value = 1651
$ws.Range("I97").Value = 1460
$ws.Range("J97").Value = 1956.6
$ws.Range("K97").Value = 1460
$ws.Range("L97").Value = 1956.6
$ws.Range("M97").Value = -964
$ws.Range("N97").Value = -2948.6

# Row 132
$ws.Range("H132").Value = 7441.0835
$ws.Range("I132").Value = 10018.267
$ws.Range("J132").Value = 3145.7778
$ws.Range("K132").Value = 30054.801
$ws.Range("L132").Value = 9437.3334
$ws.Range("M132").Value = -27524.801
$ws.Range("N132").Value = -14497.3334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 11780.167
$ws.Range("I40").Value = 13895
$ws.Range("J40").Value = 7550.5
$ws.Range("K40").Value = 13895
$ws.Range("L40").Value = 7550.5
$ws.Range("M40").Value = -13759
$ws.Range("N40").Value = -7822.5

# Row 119
$ws.Range("H119").Value = 48250
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 48250
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 48250
$ws.Range("N119").Value = -57926

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 30000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 30000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30346

# Row 132
$ws.Range("H132").Value = 1328.2678
$ws.Range("I132").Value = 1099.6111
$ws.Range("J132").Value = 7502
$ws.Range("K132").Value = 3298.8333
$ws.Range("L132").Value = 22506
$ws.Range("M132").Value = -768.8333000000002
$ws.Range("N132").Value = -27566

# Row 133
$ws.Range("H133").Value = 33782.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 33782.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 33782.5
$ws.Range("N133").Value = -43902.5
